$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 / 13: TRON and Chainlink swap positions, with updated price/volume data
# Row 2
$dCell = $ws.Cells.Item(2, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '52.371.58'
$dCell.Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  +0.25%  '

# Row 3
$dCell = $ws.Cells.Item(3, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '2.843.94'
$dCell.Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  +1.24%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.03%  '

# Row 5
$dCell = $ws.Cells.Item(5, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '362.12'
$dCell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +4.06%  '

# Row 6
$dCell = $ws.Cells.Item(6, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '112.75'
$dCell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -2.64%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  +4.05%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  -0.06%  '

# Row 9
$dCell = $ws.Cells.Item(9, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.609'
$dCell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +3.10%  '

# Row 10
$dCell = $ws.Cells.Item(10, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '41.18'
$dCell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -2.83%  '

# Row 11
$dCell = $ws.Cells.Item(11, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.0868'
$dCell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +0.61%  '

# Row 12
$ws.Cells.Item(12, 2).Value = 'Chainlink'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$dCell = $ws.Cells.Item(12, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '20.15'
$dCell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +0.56%  '

# Row 13
$ws.Cells.Item(13, 2).Value = 'TRON'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$dCell = $ws.Cells.Item(13, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.132'
$dCell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +1.09%  '

# Row 14
$dCell = $ws.Cells.Item(14, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '7.84'
$dCell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  -0.13%  '

# Row 15
$dCell = $ws.Cells.Item(15, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '3.286.70'
$dCell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +1.28%  '

# Row 16
$dCell = $ws.Cells.Item(16, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '2.804.54'
$dCell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +0.34%  '

# Row 17
$dCell = $ws.Cells.Item(17, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.937'
$dCell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +5.07%  '

# Row 18
$dCell = $ws.Cells.Item(18, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '52.317.07'
$dCell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  +0.21%  '

# Row 19
$dCell = $ws.Cells.Item(19, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '7.60'
$dCell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +3.74%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  -0.76%  '

# Row 21
$dCell = $ws.Cells.Item(21, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '13.49'
$dCell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +0.99%  '

# Row 22
$ws.Cells.Item(22, 5).Value = '  +1.97%  '

# Row 23
$dCell = $ws.Cells.Item(23, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '273.35'
$dCell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +1.23%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  +0.72%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  +2.76%  '

# Row 26
$dCell = $ws.Cells.Item(26, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '27.03'
$dCell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +0.48%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  -0.01%  '

# Row 28
$dCell = $ws.Cells.Item(28, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '10.37'
$dCell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +0.82%  '

# Row 29
$dCell = $ws.Cells.Item(29, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '2.25'
$dCell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -0.03%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  +2.23%  '

# Row 31
$dCell = $ws.Cells.Item(31, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.0484'
$dCell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +8.90%  '

# Row 32
$dCell = $ws.Cells.Item(32, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '35.48'
$dCell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +2.89%  '

# Row 33
$dCell = $ws.Cells.Item(33, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '52.59'
$dCell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +4.39%  '

# Row 34
$dCell = $ws.Cells.Item(34, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '5.91'
$dCell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +2.15%  '

# Row 35
$dCell = $ws.Cells.Item(35, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '5.57'
$dCell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +12.77%  '

# Row 36
$dCell = $ws.Cells.Item(36, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.0855'
$dCell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +2.40%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  +0.04%  '

# Row 38
$dCell = $ws.Cells.Item(38, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '3.30'
$dCell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +1.87%  '

# Row 39
$dCell = $ws.Cells.Item(39, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '2.05'
$dCell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -2.86%  '

# Row 40
$dCell = $ws.Cells.Item(40, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '18.51'
$dCell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -1.05%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  +1.85%  '

# Row 42
$dCell = $ws.Cells.Item(42, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '2.57'
$dCell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -0.99%  '

# Row 43
$dCell = $ws.Cells.Item(43, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '126.56'
$dCell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +0.13%  '

# Row 44
$dCell = $ws.Cells.Item(44, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '23.08'
$dCell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -1.43%  '

# Row 45
$dCell = $ws.Cells.Item(45, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '2.28'
$dCell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  -0.43%  '

# Row 46
$dCell = $ws.Cells.Item(46, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '2.099.31'
$dCell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +1.84%  '

# Row 47
$dCell = $ws.Cells.Item(47, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '3.36'
$dCell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +1.23%  '

# Row 48
$dCell = $ws.Cells.Item(48, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '2.30'
$dCell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -1.70%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  +5.95%  '

# Row 50
$dCell = $ws.Cells.Item(50, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.967'
$dCell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -0.01%  '

# Row 51
$dCell = $ws.Cells.Item(51, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '9.26'
$dCell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +2.97%  '
